# Apply the diff to comparison_summary.docx via Word COM-interop Find/Replace.
$d = $word.ActiveDocument

function Replace-One([string]$old, [string]$new) {
    # wdReplace:=1 (wdReplaceOne) replaces only the first match found,
    # scanning from the very beginning of the supplied range each call.
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}

# --- Exception Counts summary paragraph ---
Replace-One "  Riverside only: 365" "  Riverside only: 222"
Replace-One "  Total Riverside exceptions: 5928" "  Total Riverside exceptions: 5785"

# --- Rule Success Rate Comparison table (Riverside columns) ---
Replace-One "4957 / 5569" "4826 / 5422"
Replace-One "609 / 680" "582 / 649"
Replace-One "89.6%" "89.7%"
Replace-One "767 / 4436" "739 / 4362"
Replace-One "17.3%" "16.9%"
Replace-One "984 / 1642" "972 / 1624"
Replace-One "260 / 331" "254 / 323"
Replace-One "78.5%" "78.6%"
Replace-One "730 / 774" "674 / 714"
Replace-One "94.3%" "94.4%"
Replace-One "267 / 1136" "241 / 1044"
Replace-One "23.5%" "23.1%"

# --- Oxford-Only Exceptions by Type: swap order of the two "...: 6" lines ---
# Use temp placeholders so the two single-token swaps don't collide.
Replace-One "Strong participle must end in -en or -e: 6" "__TMP_SWAP_A__"
Replace-One "Weak pt sg must end in -ed, -d, or -t: 6" "Strong participle must end in -en or -e: 6"
Replace-One "__TMP_SWAP_A__" "Weak pt sg must end in -ed, -d, or -t: 6"

# --- Riverside-Only Exceptions by Type: values updated, order rearranged ---
Replace-One "Weak pt sg must end in -ed, -d, or -t: 245" "Weak pt sg must end in -ed, -d, or -t: 179"
Replace-One "Present 3rd sg must end in -eth: 58" "Strong pt sg must not end in -en or -e: 15"
Replace-One "Infinitive must end in -en or -e: 24" "Present 3rd sg must end in -eth: 13"
Replace-One "Strong pt sg must not end in -en or -e: 19" "Infinitive must end in -en or -e: 8"
Replace-One "Present plural must end in -en or -e: 10" "Present plural must end in -en or -e: 4"
Replace-One "Past plural must end in -en or -e: 6" "Past plural must end in -en or -e: 2"
Replace-One "Strong participle must end in -en or -e: 3" "Strong participle must end in -en or -e: 1"
